$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.884.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.536.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.29"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.12"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.59"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.924.44"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.542.71"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.927.23"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.12"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.83"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0966"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.09"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.68"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.88"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.45"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.95"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.88"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.35"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.07"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0784"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("B40").Value = "ApeXProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.31"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +13.43%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.05"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.77%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0305"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.037.22"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.16"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.18"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.41"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.74"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.777.59"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.61%  "
